$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Phone (D2) and password (E2) become plain numbers instead of text
$ws.Range("D2").Value = 5198176511
$ws.Range("E2").Value = 1234

# Hosted games (G2) gains an extra entry "1;"
$ws.Range("G2").Value = "0;1;"
